$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq "Miss Dina Nasr, Administrator") {
        $cell.Value = "Administrator, Miss Dina Nasr"
    }
}
